# The sheet had a stray leading "A" column (a duplicate index column) and a
# mislabelled header "MODEL_CONDITION". This edit:
#   1. Fixes the header text (MODEL_CONDITION -> MODELCONDITION).
#   2. Removes the stray column A entirely, shifting B:F left to A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the mislabelled header text in place.
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# 2. Delete the stray leading column; remaining columns shift left
#    (old B:F -> new A:E), matching the target layout exactly.
$ws.Columns("A:A").Delete()
